$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.676.85"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "2.293.66"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "96.47"
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "268.96"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("E9").Value = "  -2.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.25"
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0928"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.88"
$ws.Range("E12").Value = "  -3.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.64"
$ws.Range("E14").Value = "  +1.47%  "
$ws.Range("D15").Value = "2.638.89"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.853"
$ws.Range("E16").Value = "  -1.39%  "
$ws.Range("D17").Value = "2.292.01"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "43.714.92"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000112"
$ws.Range("E19").Value = "  +3.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.19"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.96"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("E22").Value = "  +10.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.24"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.10"
$ws.Range("E24").Value = "  -5.45%  "
$ws.Range("E25").Value = "  +5.45%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.23"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.55"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.87"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.81"
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0904"
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.40"
$ws.Range("E34").Value = "  -1.99%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.53"
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0351"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.106"
$ws.Range("E38").Value = "  -2.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.36"
$ws.Range("E39").Value = "  -2.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.237"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.34"
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.11"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.49"
$ws.Range("E44").Value = "  +4.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.78"
$ws.Range("E45").Value = "  -2.87%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.102"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.16"
$ws.Range("E47").Value = "  -5.34%  "
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.05"
$ws.Range("E49").Value = "  -3.48%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.444"
$ws.Range("E50").Value = "  +5.32%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.53"
$ws.Range("E51").Value = "  +12.87%  "
